$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell A234 ("taweno" -> "taweno | ta bueno")
$ws.Range('A234').Value = 'taweno | ta bueno'

# Add new meme/stickerID rows (235-265)
$ws.Range('A235').Value = 'fino señores | fino'
$ws.Range('B235').Value = 'CAACAgEAAxkBAAKYnmJ6iEQx16bDPiv_Zk-T8eErj5PfAAI2AwACz5TRRwQWtGItC9rzJAQ'
$ws.Range('A236').Value = 'un clasico | clasico'
$ws.Range('B236').Value = 'CAACAgEAAxkBAAKYoGJ6iFKC6gM1A24egnKmtTGGjbWSAAKGAgACXOnZR9ugI6YifSYPJAQ'
$ws.Range('A237').Value = 'Por que no me muero'
$ws.Range('B237').Value = 'CAACAgEAAxkBAAKYomJ6iHLn8KVYVkD0PA2AC2ODAAGilQACrwEAAmzl0EfznyZBVvMFayQE'
$ws.Range('A238').Value = 'siempre espera lo inesperado | espera lo inesperado'
$ws.Range('B238').Value = 'CAACAgEAAxkBAAKYpGJ6iIkE9L8eovhds9nrpwWsx4ylAALHAQACD9_YR3fm0WEq6sAVJAQ'
$ws.Range('A239').Value = 'algo anda mal'
$ws.Range('B239').Value = 'CAACAgEAAxkBAAKYqGJ6iPlPiufwmlO3zsr7eLwKyQqvAAKOAgAC0LjQRwU7oeZ8Xt4nJAQ'
$ws.Range('A240').Value = 'estuviste practicando | estuviste practicando eh'
$ws.Range('B240').Value = 'CAACAgEAAxkBAAKYqmJ6iQXX05Kno-OX60bgrUn6ftjKAAJUAwAC88TRR967aBO56HNIJAQ'
$ws.Range('A241').Value = 'ya no quiero verte nunca mas'
$ws.Range('B241').Value = 'CAACAgEAAxkBAAKYrGJ6iRhnrFVQ3s0Glku_uqGPRgSGAAIZAgAC5I_RR0d9A0gIo4PNJAQ'
$ws.Range('A242').Value = 'bien pensado woody | bien pensado'
$ws.Range('B242').Value = 'CAACAgEAAxkBAAKYrmJ6iSiD-Ph0eXxQhrYP-ZW4EhKNAAJCAgACDiTZR0aQ52guLvymJAQ'
$ws.Range('A243').Value = 'debe haber un mensaje oculto en algun lado'
$ws.Range('B243').Value = 'CAACAgEAAxkBAAKYsGJ6iTtPhFZ8kPB8CF39hJLBjhB4AAL3BQACNKLQR7QiSv-s8bKgJAQ'
$ws.Range('A244').Value = 'adinivare escuela publica | escuela publica?'
$ws.Range('B244').Value = 'CAACAgEAAxkBAAKYsmJ6iVA2ZrRZS0YC1g7K7ZdD-yicAAJqAgACYiHRRybYHd7jZ7XiJAQ'
$ws.Range('A245').Value = 'que pendejo | pendejo | bien pendejo'
$ws.Range('B245').Value = 'CAACAgEAAxkBAAKYtGJ6iWf2dgWtNFaeGhLD5yiwiuwDAAIHAgACqgXYR_F5M_vV5nRFJAQ'
$ws.Range('A246').Value = 'yo no recibo ordenes | yo no recibo ordenes de ti | soy un prime yo no recibo ordenes de ti'
$ws.Range('B246').Value = 'CAACAgEAAxkBAAKYtmJ6iYB_6WUTDu8EGLRsc3dXqYA5AAIFAwAC-jXYRzbXRJG6bBIwJAQ'
$ws.Range('A247').Value = 'chingues a tu madre'
$ws.Range('B247').Value = 'CAACAgEAAxkBAAKYuGJ6iaf1-S-q5R696tl3pepMYXznAAKxAgACSE3ZR1Jp-2qiDk2QJAQ'
$ws.Range('A248').Value = 'las cosas se salieron de control'
$ws.Range('B248').Value = 'CAACAgEAAxkBAAKYumJ6ibVWLKDWhicUbr0LSRgvOSVkAAJ3AgAC0CDZR3KDm6bifm8lJAQ'
$ws.Range('A249').Value = 'norman se fue de sabatico | se fue se sabatico | norman se fue de sabatico cariño'
$ws.Range('B249').Value = 'CAACAgEAAxkBAAKYvGJ6icMFICUY1uR0I_hVAkrnTMM3AAKFAgACwfXRRywaOzhLEx9xJAQ'
$ws.Range('A250').Value = 'apagalo otto | apagalo otto apagalo'
$ws.Range('B250').Value = 'CAACAgEAAxkBAAKYvmJ6idoPsC3Tr5MJnqGJj5cDTskbAALVAQAC9GzRR7ZV9S4QakdLJAQ'
$ws.Range('A251').Value = 'te convertiste en aquello que juraste destruir'
$ws.Range('B251').Value = 'CAACAgEAAxkBAAKYwGJ6ieqYwhLW9l8QqNSF0OsVriteAALYAQACXGHYRxvjCFmihjrXJAQ'
$ws.Range('A252').Value = 'acompañame a ver esta triste historia | esta triste historia'
$ws.Range('B252').Value = 'CAACAgEAAxkBAAKYwmJ6iiH05krRZekYLrLCsaqzyH3qAAIPAgAC9tfZR1jHe_5ngytbJAQ'
$ws.Range('A253').Value = 'ya no hijo no te creo | no te creo'
$ws.Range('B253').Value = 'CAACAgEAAxkBAAKYxGJ6ikvscJegPShuDH5mSFaYQcnxAAIbAgACUeXQR499OJhOSkwAASQE'
$ws.Range('A254').Value = 'primera vez'
$ws.Range('B254').Value = 'CAACAgEAAxkBAAKYymJ6i9wAAbYjMRv--cB8wI9Z4iT_HAACHAIAAjr32EdYdqbqqhF0RSQE'
$ws.Range('A255').Value = 'cuanta virgindad hay en este mundo'
$ws.Range('B255').Value = 'CAACAgEAAxkBAAKYzGJ6i-vNBvCeeQ3oCGgTbWJRYkYHAAJjAwACkhzRR2z7JrhL1e33JAQ'
$ws.Range('A256').Value = 'nel'
$ws.Range('B256').Value = 'CAACAgEAAxkBAAKaumKDSUnRGEQjfMQdX682-ByY63sOAAJxAgAC284YRKLIe-5NKoKYJAQ'
$ws.Range('A257').Value = 'yo te conozco'
$ws.Range('B257').Value = 'CAACAgEAAxkBAAKavGKDSXmOurwJrmwaiGhtOzZbm-0sAAJAAwACczoYROdUjUNYoprbJAQ'
$ws.Range('A258').Value = 'que se armen los pinches chingadazos | que se armen'
$ws.Range('B258').Value = 'CAACAgEAAxkBAAKavmKDSYavYJfvqAIqDs-l1FcRoGxOAAJMAgAC5zcZRPOHQ_snEfetJAQ'
$ws.Range('A259').Value = 'corre perra corre'
$ws.Range('B259').Value = 'CAACAgEAAxkBAAKawGKDSaa8u8rh8QeCgweeW8SxawQkAAIfAgAC6cMYRFCinDUjQcOWJAQ'
$ws.Range('A260').Value = 'lo voy a gozar'
$ws.Range('B260').Value = 'CAACAgEAAxkBAAKawmKDSbgcqCoo5iMYYuaViyh4DgU0AALFAgACrgUZRM1F_oBmJYN1JAQ'
$ws.Range('A261').Value = 'besitos besitos chau chau'
$ws.Range('B261').Value = 'CAACAgEAAxkBAAKaxGKDScWmNTZ5OJKhF7inUqLKqPW5AALjAgACfG4YRDBqUQrUUmSWJAQ'
$ws.Range('A262').Value = 'le falle señor | le falle | le falle señor ud confió en mi y yo le falle | ud confio en mi y yo le falle'
$ws.Range('B262').Value = 'CAACAgEAAxkBAAKaxmKDSdQguJfFnnK9Yz8-VqWkxzKRAAINAwAC37wYRDwz3FbdcjwSJAQ'
$ws.Range('A263').Value = 'valio madres | ya valio madres'
$ws.Range('B263').Value = 'CAACAgEAAxkBAAKayGKDSfljGeVupdp6TF1L4Rkudco_AAKDAwACqyEZRLgBZeV5MIg2JAQ'
$ws.Range('A264').Value = 'obvio'
$ws.Range('B264').Value = 'CAACAgEAAxkBAAKaymKDSgkznlteVbYTrY6BghO5eO4fAAK3AQAC6iwhRB0MhmFrlUETJAQ'
$ws.Range('A265').Value = 'unlimited power'
$ws.Range('B265').Value = 'CAACAgEAAxkBAAKazGKDShMnDI8rm45-pE_AHu1jiFz4AALNAgACOKEZRJyc5mNdrGO-JAQ'
